$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Charts")

# Update the "Story points" (B) values for sprints 0-3.
# B6 currently holds the formula =900*3 (2700); replace with a plain value.
$ws.Range("B6").Value = 2350
$ws.Range("B7").Value = 1000
$ws.Range("B8").Value = 1000
$ws.Range("B9").Value = 1000

# Move the active selection to B7 (was B10).
$ws.Activate()
$ws.Range("B7").Select()
